$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 10000
$ws.Range("J7").Value = 10000
$ws.Range("L7").Value = 10000
$ws.Range("N7").Value = -10224
$ws.Range("H13").Value = 100006
$ws.Range("J13").Value = 100006
$ws.Range("L13").Value = 100006
$ws.Range("N13").Value = -100344
$ws.Range("H14").Value = 10000
$ws.Range("J14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("N14").Value = -10382
$ws.Range("H16").Value = 70006.664
$ws.Range("J16").Value = 70006.664
$ws.Range("L16").Value = 70006.664
$ws.Range("N16").Value = -70466.664
$ws.Range("H17").Value = 4210.3887
$ws.Range("J17").Value = 4210.3887
$ws.Range("L17").Value = 12631.1661
$ws.Range("N17").Value = -12967.1661
$ws.Range("H18").Value = 2315839.2
$ws.Range("I18").Value = 2525979.2
$ws.Range("K18").Value = 2525979.2
$ws.Range("M18").Value = -2525695.2
$ws.Range("H19").Value = 6211674.5
$ws.Range("I19").Value = 12987461
$ws.Range("J19").Value = 536.75
$ws.Range("K19").Value = 12987461
$ws.Range("L19").Value = 536.75
$ws.Range("M19").Value = -12987286
$ws.Range("N19").Value = -886.75
$ws.Range("H20").Value = 39994.832
$ws.Range("I20").Value = 9980.25
$ws.Range("J20").Value = 100024
$ws.Range("K20").Value = 9980.25
$ws.Range("L20").Value = 100024
$ws.Range("M20").Value = -9750.25
$ws.Range("N20").Value = -100484
$ws.Range("H21").Value = 28339.666
$ws.Range("I21").Value = 50019
$ws.Range("J21").Value = 17500
$ws.Range("K21").Value = 50019
$ws.Range("L21").Value = 17500
$ws.Range("M21").Value = -49551
$ws.Range("N21").Value = -18436
$ws.Range("H23").Value = 28339.666
$ws.Range("I23").Value = 50019
$ws.Range("J23").Value = 17500
$ws.Range("K23").Value = 50019
$ws.Range("L23").Value = 17500
$ws.Range("M23").Value = -49785
$ws.Range("N23").Value = -17968
$ws.Range("H26").Value = 61346
$ws.Range("J26").Value = 61346
$ws.Range("L26").Value = 61346
$ws.Range("N26").Value = -62034
$ws.Range("H29").Value = 5000
$ws.Range("J29").Value = 5000
$ws.Range("L29").Value = 15000
$ws.Range("N29").Value = -15562
$ws.Range("H31").Value = 900.5
$ws.Range("I31").Value = 900.5
$ws.Range("K31").Value = 2701.5
$ws.Range("M31").Value = -2471.5
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H33").Value = 839.1905
$ws.Range("I33").Value = 667.93335
$ws.Range("K33").Value = 667.93335
$ws.Range("M33").Value = -438.93335
$ws.Range("H34").Value = 16537.23
$ws.Range("I34").Value = 1531
$ws.Range("J34").Value = 29399.715
$ws.Range("K34").Value = 1531
$ws.Range("L34").Value = 29399.715
$ws.Range("M34").Value = -1328
$ws.Range("N34").Value = -29805.715
$ws.Range("H35").Value = 39994.832
$ws.Range("I35").Value = 9980.25
$ws.Range("J35").Value = 100024
$ws.Range("K35").Value = 9980.25
$ws.Range("L35").Value = 100024
$ws.Range("M35").Value = -9601.25
$ws.Range("N35").Value = -100782
$ws.Range("H36").Value = 16537.23
$ws.Range("I36").Value = 1531
$ws.Range("J36").Value = 29399.715
$ws.Range("K36").Value = 1531
$ws.Range("L36").Value = 29399.715
$ws.Range("M36").Value = -816
$ws.Range("N36").Value = -30829.715
$ws.Range("H132").Value = 6456033.5
$ws.Range("I132").Value = 8336877.5
$ws.Range("J132").Value = 7425.5713
$ws.Range("K132").Value = 25010632.5
$ws.Range("L132").Value = 22276.7139
$ws.Range("M132").Value = -25008102.5
$ws.Range("N132").Value = -27336.7139

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1311.8379
$ws.Range("I110").Value = 526
$ws.Range("J110").Value = 3169.2727
$ws.Range("K110").Value = 526
$ws.Range("L110").Value = 3169.2727
$ws.Range("M110").Value = 1519
$ws.Range("N110").Value = -7259.2727

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3702.7778
$ws.Range("I134").Value = 3332.353
$ws.Range("K134").Value = 9997.059000000001
$ws.Range("M134").Value = -7462.059000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1780.07
$ws.Range("I31").Value = 1053.4043
$ws.Range("J31").Value = 2424.4717
$ws.Range("K31").Value = 1053.4043
$ws.Range("L31").Value = 2424.4717
$ws.Range("M31").Value = -758.4042999999999
$ws.Range("N31").Value = -3014.4717
$ws.Range("H34").Value = 1780.07
$ws.Range("I34").Value = 1053.4043
$ws.Range("J34").Value = 2424.4717
$ws.Range("K34").Value = 1053.4043
$ws.Range("L34").Value = 2424.4717
$ws.Range("M34").Value = -851.4042999999999
$ws.Range("N34").Value = -2828.4717
$ws.Range("H58").Value = 10873007
$ws.Range("I58").Value = 2679.9546
$ws.Range("J58").Value = 20837474
$ws.Range("K58").Value = 2679.9546
$ws.Range("L58").Value = 20837474
$ws.Range("M58").Value = -2476.9546
$ws.Range("N58").Value = -20837880
$ws.Range("H132").Value = 3245.875
$ws.Range("I132").Value = 3077.7856
$ws.Range("J132").Value = 3481.2
$ws.Range("K132").Value = 9233.356800000001
$ws.Range("L132").Value = 10443.6
$ws.Range("M132").Value = -6703.356800000001
$ws.Range("N132").Value = -15503.6
$ws.Range("H134").Value = 22731166
$ws.Range("I134").Value = 31252928
$ws.Range("J134").Value = 6466.5
$ws.Range("K134").Value = 93758784
$ws.Range("L134").Value = 19399.5
$ws.Range("M134").Value = -93756249
$ws.Range("N134").Value = -24469.5
$ws.Range("H136").Value = 10873007
$ws.Range("I136").Value = 2679.9546
$ws.Range("J136").Value = 20837474
$ws.Range("K136").Value = 8039.8638
$ws.Range("L136").Value = 62512422
$ws.Range("M136").Value = -5489.8638
$ws.Range("N136").Value = -62517522

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2576.1462
$ws.Range("I132").Value = 1697.7693
$ws.Range("K132").Value = 5093.3079
$ws.Range("M132").Value = -2563.3079
$ws.Range("H136").Value = 3452774.2
$ws.Range("I136").Value = 5886144
$ws.Range("J136").Value = 5500.4165
$ws.Range("K136").Value = 17658432
$ws.Range("L136").Value = 16501.2495
$ws.Range("M136").Value = -17655882
$ws.Range("N136").Value = -21601.2495

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1890.909
$ws.Range("J81").Value = 2250
$ws.Range("L81").Value = 4500
$ws.Range("N81").Value = -6622
$ws.Range("H84").Value = 1890.909
$ws.Range("J84").Value = 2250
$ws.Range("L84").Value = 22500
$ws.Range("N84").Value = -33108
